$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.851.64'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '3.111.98'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '3.108.96'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('E10').Value = '  -3.68%  '
$ws.Range('E11').Value = '  -1.93%  '
$ws.Range('E13').Value = '  -2.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.36'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('E15').Value = '  -1.45%  '
$ws.Range('D16').Value = '3.627.59'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').Value = '66.793.13'
$ws.Range('E17').Value = '  -0.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.13'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = '3.111.48'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('E20').Value = '  +2.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '477.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.714'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.56'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.91'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.41%  '
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.61'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  -1.74%  '
$ws.Range('D34').Value = '0.0₃0941'
$ws.Range('E34').Value = '  -8.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.85'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.973'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '46.99'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.88%  '
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.309'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.90%  '
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.67'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('D44').Value = '2.810.81'
$ws.Range('E44').Value = '  +1.35%  '
$ws.Range('E45').Value = '  -2.36%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '380.15'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.55'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -11.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '136.28'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.98'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('E51').Value = '  -2.04%  '
